$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing H-column values in the "Growth" column ---

# H20: plain 0.5 -> formula =1/3
$ws.Range("H20").Formula = "=1/3"

# H23: 2.5 -> 0.7
$ws.Range("H23").Value = 0.7

# H24: 0.5 -> 0.9
$ws.Range("H24").Value = 0.9

# H25: 0.5 -> 1
$ws.Range("H25").Value = 1

# H26: formula =2/6 -> plain value 1.5
$ws.Range("H26").Value = 1.5

# H27: 0.25 -> formula =4/3
$ws.Range("H27").Formula = "=4/3"

# H28: formula =2/3 -> plain value 2
$ws.Range("H28").Value = 2

# H29: formula =7/3 -> plain value 3
$ws.Range("H29").Value = 3

# H126: 0 -> formula =2/3
$ws.Range("H126").Formula = "=2/3"

# H168: 1.35 -> 0.35
$ws.Range("H168").Value = 0.35

# H169: 2.5 -> 0.5
$ws.Range("H169").Value = 0.5

# H216: 2.5299999999999998 -> 1.53
$ws.Range("H216").Value = 1.53

# --- Append 7 new rows of data (rows 618-624) ---
# Each row is a hashtable keyed by column letter (column F is always a
# shared ABS(D-E) formula, so it is not listed here and is handled below).

$newRows = @(
    @{ A=45874; B="Flowering";    C="Large";  D=65; E=83; G=0; H=0;   I="No"; J=2; K="Bright";  L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Nonflowering"; C="Medium"; D=65; E=83; G=0; H=0;   I="No"; J=3; K="Bright";  L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Nonflowering"; C="Small";  D=65; E=83; G=0; H=0;   I="No"; J=3; K="Bright";  L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Nonflowering"; C="Medium"; D=65; E=83; G=0; H=0;   I="No"; J=3; K="Neutral"; L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Nonflowering"; C="Medium"; D=65; E=83; G=0; H=0.1; I="No"; J=3; K="Neutral"; L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Nonflowering"; C="Large";  D=65; E=83; G=0; H=0.1; I="No"; J=4; K="Neutral"; L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 },
    @{ A=45874; B="Tree";         C="Medium"; D=65; E=83; G=0; H=0.5; I="No"; J=1; K="Neutral"; L=7; M=0.46; N=60; O=30.3; P=16; Q=0.13; R=8.6999999999999993; S=75; T=4 }
)

$startRow = 618
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}

# Column F ("Temp_Diff") continues the existing ABS(D-E) pattern used for
# every other row in the sheet. Assigning the relative formula across the
# whole new block in one go (rather than cell-by-cell) lets it land as a
# single shared-formula group, mirroring how the pre-existing F607:F617
# block (and all the other Fxx:Fyy blocks above it) are stored.
$ws.Range("F" + $startRow + ":F" + ($startRow + $newRows.Count - 1)).Formula = "=ABS(D" + $startRow + "-E" + $startRow + ")"

# --- Update the view so the new rows are visible ---
$ws.Range("A599").Select()
$ws.Range("I625").Select()
